# Update column F ("想去人数" / interest counts) per commit 456a3b4
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 342   # was 341
$ws.Range("F3").Value = 231   # was 229
$ws.Range("F4").Value = 558   # was 556
$ws.Range("F5").Value = 1342   # was 1339
$ws.Range("F6").Value = 653   # was 652
$ws.Range("F7").Value = 348   # was 347
$ws.Range("F11").Value = 6196   # was 6185
$ws.Range("F12").Value = 114   # was 113
$ws.Range("F13").Value = 23   # was 21
$ws.Range("F14").Value = 1892   # was 1891
$ws.Range("F15").Value = 4640   # was 4634
$ws.Range("F18").Value = 305   # was 306
$ws.Range("F19").Value = 5411   # was 5404
$ws.Range("F20").Value = 7068   # was 7054
$ws.Range("F22").Value = 1085   # was 1086
$ws.Range("F23").Value = 750   # was 749
$ws.Range("F24").Value = 3978   # was 3973
$ws.Range("F25").Value = 550   # was 548
$ws.Range("F29").Value = 1051   # was 1050
$ws.Range("F30").Value = 1490   # was 1487
$ws.Range("F32").Value = 681   # was 680
$ws.Range("F33").Value = 1679   # was 1675
$ws.Range("F35").Value = 1877   # was 1871
$ws.Range("F36").Value = 228   # was 226
$ws.Range("F38").Value = 1235   # was 1232
$ws.Range("F40").Value = 680   # was 677
$ws.Range("F42").Value = 909   # was 845
$ws.Range("F43").Value = 3648   # was 3643
$ws.Range("F48").Value = 88   # was 87

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1262   # was 1260
$ws.Range("F5").Value = 44   # was 43

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4366   # was 4360

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4366   # was 4360
$ws.Range("F3").Value = 342   # was 341
$ws.Range("F4").Value = 1262   # was 1260
$ws.Range("F6").Value = 231   # was 229
$ws.Range("F7").Value = 558   # was 556
$ws.Range("F9").Value = 1342   # was 1339
$ws.Range("F11").Value = 653   # was 652
$ws.Range("F12").Value = 348   # was 347
$ws.Range("F15").Value = 114   # was 113
$ws.Range("F16").Value = 23   # was 21
$ws.Range("F17").Value = 4640   # was 4634
$ws.Range("F18").Value = 5412   # was 5404
$ws.Range("F19").Value = 5412   # was 5404
$ws.Range("F21").Value = 1085   # was 1086
$ws.Range("F22").Value = 750   # was 749
$ws.Range("F23").Value = 3978   # was 3973
$ws.Range("F24").Value = 550   # was 548
$ws.Range("F28").Value = 1051   # was 1050
$ws.Range("F29").Value = 1490   # was 1487
$ws.Range("F31").Value = 681   # was 680
$ws.Range("F32").Value = 1679   # was 1675
$ws.Range("F34").Value = 1877   # was 1871
$ws.Range("F39").Value = 680   # was 677
$ws.Range("F43").Value = 3648   # was 3643
$ws.Range("F48").Value = 88   # was 87
